$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "29.273.62"
$ws.Range("E2").Value = "  +2.73%  "

# Row 3
$ws.Range("D3").Value = "1.902.05"
$ws.Range("E3").Value = "  +1.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -1.58%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.45"
$ws.Range("E5").Value = "  -0.80%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  -1.66%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5139"
$ws.Range("E7").Value = "  +0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3935"
$ws.Range("E8").Value = "  -0.77%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08453"
$ws.Range("E9").Value = "  -0.45%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.53"
$ws.Range("E10").Value = "  +1.11%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.119"
$ws.Range("E11").Value = "  +0.74%  "

# Row 12
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.904.12"
$ws.Range("E12").Value = "  +1.66%  "

# Row 13
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.252"
$ws.Range("E13").Value = "  -0.12%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.66"
$ws.Range("E14").Value = "  +0.79%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.362"
$ws.Range("E15").Value = "  +1.49%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.005"
$ws.Range("E16").Value = "  -1.61%  "

# Row 17
$ws.Range("E17").Value = "  +2.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001107"
$ws.Range("E18").Value = "  -0.51%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06731"
$ws.Range("E19").Value = "  -1.01%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.91"
$ws.Range("E20").Value = "  +1.03%  "

# Row 21
$ws.Range("E21").Value = "  -1.54%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.032"
$ws.Range("E22").Value = "  +1.48%  "

# Row 23
$ws.Range("D23").Value = "29.266.57"
$ws.Range("E23").Value = "  +2.57%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.12"
$ws.Range("E24").Value = "  -0.62%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.220"
$ws.Range("E25").Value = "  -3.29%  "

# Row 26
$ws.Range("D26").Value = "2.120.23"
$ws.Range("E26").Value = "  +1.57%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "160.28"
$ws.Range("E27").Value = "  -1.10%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.99"
$ws.Range("E28").Value = "  +0.86%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.438"
$ws.Range("E29").Value = "  +3.58%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "127.65"
$ws.Range("E30").Value = "  +0.38%  "

# Row 31
$ws.Range("B31").Value = "Stellar"
$ws.Range("C31").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1050"
$ws.Range("E31").Value = "  -0.71%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.058"
$ws.Range("E32").Value = "  +1.17%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.006"
$ws.Range("E33").Value = "  +3.97%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.650"
$ws.Range("E34").Value = "  -0.06%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.02478"
$ws.Range("E35").Value = "  +1.75%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06587"
$ws.Range("E36").Value = "  +1.48%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.137"
$ws.Range("E37").Value = "  +2.87%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2195"
$ws.Range("E38").Value = "  +0.70%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.235"
$ws.Range("E39").Value = "  +4.12%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.121"
$ws.Range("E40").Value = "  +2.29%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6526"
$ws.Range("E41").Value = "  +1.95%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.238"
$ws.Range("E42").Value = "  -2.40%  "

# Row 43
$ws.Range("E43").Value = "  +0.42%  "

# Row 44
$ws.Range("B44").Value = "Decentraland"
$ws.Range("C44").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6062"
$ws.Range("E44").Value = "  +0.11%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.23"
$ws.Range("E45").Value = "  +1.67%  "

# Row 46
$ws.Range("B46").Value = "PancakeSwap"
$ws.Range("C46").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.677"
$ws.Range("E46").Value = "  -1.36%  "

# Row 47
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.055"
$ws.Range("E47").Value = "  +2.99%  "

# Row 48
$ws.Range("B48").Value = "EOS"
$ws.Range("C48").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.232"
$ws.Range("E48").Value = "  +2.04%  "

# Row 49
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.03"
$ws.Range("E49").Value = "  +0.58%  "

# Row 50
$ws.Range("B50").Value = "WEMIXTOKEN"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.177"
$ws.Range("E50").Value = "  -2.56%  "

# Row 51
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "77.86"
$ws.Range("E51").Value = "  +1.53%  "
